$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: 2022 data, appended after the existing 2021 (column K) data.

# Header row (row 4) - year label, copy style from K4
$ws.Range("L4").Value = 2022
$ws.Range("L4").Style = $ws.Range("K4").Style

# Row 5 - total withdrawal
$ws.Range("L5").Value = 8800.6
$ws.Range("L5").Style = $ws.Range("K5").Style

# Row 6 - section header (by type of source) - empty, copy style only
$ws.Range("L6").Style = $ws.Range("K6").Style

# Row 7 - from natural water sources (formula)
$ws.Range("L7").Formula = "=L5-L8"
$ws.Range("L7").Style = $ws.Range("K7").Style

# Row 8 - from underground horizons
$ws.Range("L8").Value = 258.39999999999998
$ws.Range("L8").Style = $ws.Range("K8").Style

# Row 9 - section header (by territory) - empty, copy style only
$ws.Range("L9").Style = $ws.Range("K9").Style

# Row 10 - Batken oblast
$ws.Range("L10").Value = 683.8
$ws.Range("L10").Style = $ws.Range("K10").Style

# Row 11 - Jalal-Abat oblast
$ws.Range("L11").Value = 1101.8
$ws.Range("L11").Style = $ws.Range("K11").Style

# Row 12 - Yssyk-Kul oblast
$ws.Range("L12").Value = 714.9
$ws.Range("L12").Style = $ws.Range("K12").Style

# Row 13 - Naryn oblast
$ws.Range("L13").Value = 757.9
$ws.Range("L13").Style = $ws.Range("K13").Style

# Row 14 - Osh oblast
$ws.Range("L14").Value = 1383.3
$ws.Range("L14").Style = $ws.Range("K14").Style

# Row 15 - Talas oblast
$ws.Range("L15").Value = 1023.7
$ws.Range("L15").Style = $ws.Range("K15").Style

# Row 16 - Chui oblast
$ws.Range("L16").Value = 2929.3
$ws.Range("L16").Style = $ws.Range("K16").Style

# Row 17 - Bishkek city
$ws.Range("L17").Value = 148.9
$ws.Range("L17").Style = $ws.Range("K17").Style

# Row 18 - Osh city
$ws.Range("L18").Value = 57
$ws.Range("L18").Style = $ws.Range("K18").Style

# Update selection as in the target workbook
$ws.Range("M4").Select()
